$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EEU data")

# Update the example preferences path values
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = -0.4

# Reflect the author's new selection / scroll position on the sheet
$ws.Activate()
$ws.Range("J2").Select()
